$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Sheet restructuring:
#    - rename "sensor_data" -> "CH4_sensor_data"
#    - insert a new "CO2_sensor_data" sheet right after it (same CH4 schema)
#    - insert a new "pesticide" sheet right before it (after "fertilizer")
# ---------------------------------------------------------------------------
$ch4 = $wb.Worksheets.Item("sensor_data")
$ch4.Name = "CH4_sensor_data"

$co2 = $wb.Worksheets.Add($null, $ch4)
$co2.Name = "CO2_sensor_data"

$pesticide = $wb.Worksheets.Add($ch4, $null)
$pesticide.Name = "pesticide"

# ---------------------------------------------------------------------------
# 2. CO2_sensor_data: same schema as CH4_sensor_data (ppm / time columns)
# ---------------------------------------------------------------------------
$co2.Range("A1").Value = "ppm"
$co2.Range("B1").Value = "DOUBLE"
$co2.Range("C1").Value = "偵測濃度"

$co2.Range("A2").Value = "time"
$co2.Range("B2").Value = "TIMESTAMP"
$co2.Range("C2").Value = "自動更新時間"

$co2.Columns.Item(1).ColumnWidth = 11.44140625
$co2.Columns.Item(2).ColumnWidth = 16.6640625
$co2.Columns.Item(3).ColumnWidth = 19.33203125

$co2.Range("A4").Select()

# ---------------------------------------------------------------------------
# 3. CH4_sensor_data: keep the same data, just resize columns a bit (the
#    sheet kept its own data, only its name/position/view changed)
# ---------------------------------------------------------------------------
$ch4.Columns.Item(1).ColumnWidth = 11.5546875
$ch4.Columns.Item(2).ColumnWidth = 16.5546875
$ch4.Columns.Item(3).ColumnWidth = 33.109375

$ch4.Range("D24").Select()

# ---------------------------------------------------------------------------
# 4. pesticide: name / unit / co2e schema (no N/P/K nutrient columns)
# ---------------------------------------------------------------------------
$pesticide.Range("A1").Value = "name"
$pesticide.Range("B1").Value = "VARCHAR"
$pesticide.Range("C1").Value = "農藥名稱"

$pesticide.Range("A2").Value = "unit"
$pesticide.Range("B2").Value = "VARCHAR"
$pesticide.Range("C2").Value = "單位"

$pesticide.Range("A3").Value = "co2e"
$pesticide.Range("B3").Value = "DOUBLE"
$pesticide.Range("C3").Value = "co2e"

$pesticide.Columns.Item(1).ColumnWidth = 16.5546875
$pesticide.Columns.Item(2).ColumnWidth = 15.5546875
$pesticide.Columns.Item(3).ColumnWidth = 19.21875

$pesticide.PageSetup.PaperSize = 9
$pesticide.PageSetup.Orientation = 1

$pesticide.Activate()
$excel.ActiveWindow.Zoom = 115
$pesticide.Range("C5").Select()

# ---------------------------------------------------------------------------
# 5. fertilizer: rename the "name" column description and add a co2e row
# ---------------------------------------------------------------------------
$fertilizer = $wb.Worksheets.Item("fertilizer")
$fertilizer.Range("C1").Value = "肥料名稱"

$fertilizer.Range("A6").Value = "co2e"
$fertilizer.Range("B6").Value = "DOUBLE"
$fertilizer.Range("C6").Value = "co2e"

$fertilizer.Activate()
$excel.ActiveWindow.Zoom = 145
$fertilizer.Range("A6:C6").Select()
$fertilizer.Range("C6").Activate()

# ---------------------------------------------------------------------------
# 6. product_imformation: rebuild the column-schema table with the new,
#    split fertilizer / pesticide fields.
# ---------------------------------------------------------------------------
$product = $wb.Worksheets.Item("product_imformation")

# clear previous schema rows (old table had 6 rows)
$product.Range("A1:C6").ClearContents()

$rows = @(
    @("creater",            "VARCHAR", "農友(輸入者)"),
    @("grow_crops",         "VARCHAR", "種植的農作物"),
    @("origin_place ",      "VARCHAR", "產地"),
    @("area",               "DOUBLE",  "種植面積"),
    @("fertilizer",         "VARCHAR", "使用的肥料"),
    @("dosage_fertilizer",  "DOUBLE",  "肥料的劑量"),
    @("pesticide",          "VARCHAR", "使用的農藥"),
    @("dosage_pesticide",   "DOUBLE",  "農藥的劑量"),
    @("fertilizer_co2e",    "DOUBLE",  "使用的肥料所產生的CO2E"),
    @("pesticide_co2e",     "DOUBLE",  "使用的農藥所產生的CO2E"),
    @("final_co2e",         "DOUBLE",  "總CO2E")
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 1
    $product.Range("A$r").Value = $rows[$i][0]
    $product.Range("B$r").Value = $rows[$i][1]
    $product.Range("C$r").Value = $rows[$i][2]

    if ($r -ne 3) {
        $product.Range("A$r").Font.Color = 0
    }
    if ($r -ne 1) {
        $product.Range("B$r").Font.Color = 7434609
    }
}

$product.Activate()
$excel.ActiveWindow.Zoom = 145
$product.Range("C5").Select()
